$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.735270261764526
$ws.Range("B1").Value = 3.082437753677368
$ws.Range("C1").Value = 2.804192066192627
$ws.Range("D1").Value = 1.884963274002075
$ws.Range("E1").Value = 0.8557856678962708
